# This workbook's data rows (2-30) get reshuffled: each row's
# Fecha/Variedad/Calidad/Volumen/Precio min/Precio max/Precio promedio/
# Unidad de comercializacion/Origen/Precio $/Kg/Kg por unidad values move
# to a different row (a permutation of the 29 data rows), while the
# Mercado/Region/Codreg/Tipo/Producto/Categoria columns stay put (they are
# identical for every row anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row(after) -> row(before) : which old row's data this row ends up with
$mapping = @{
    2 = 24; 3 = 17; 4 = 12; 5 = 6; 6 = 20; 7 = 11; 8 = 16; 9 = 23; 10 = 30;
    11 = 25; 12 = 19; 13 = 28; 14 = 22; 15 = 10; 16 = 8; 17 = 13; 18 = 26;
    19 = 21; 20 = 3; 21 = 4; 22 = 15; 23 = 14; 24 = 18; 25 = 5; 26 = 2;
    27 = 29; 28 = 7; 29 = 27; 30 = 9
}

# Columns whose values travel together with a row in the permutation.
$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)   # D,K,L,M,N,O,P,Q,R,S,T

# Snapshot the "before" values for every affected column/row first, since
# the mapping is a permutation (multiple interleaved cycles) and not a
# simple shift, so writes must not clobber data that is still needed.
$snapshot = @{}
for ($r = 2; $r -le 30; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
